# Updates cryptocurrency Price (D) and Volume(1h) (E) columns
# to match the latest scrape, mirroring the committed XML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.147.49"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.669.02"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "'209.57"
$ws.Range("E5").Value = "  -3.75%  "
$ws.Range("D6").Value = "'0.5227"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").Value = "'0.2616"
$ws.Range("E8").Value = "  -3.71%  "
$ws.Range("D9").Value = "'0.06337"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'21.14"
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("D11").Value = "'0.07531"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("D12").Value = "1.672.95"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "'4.433"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").Value = "'0.5503"
$ws.Range("E14").Value = "  -4.92%  "
$ws.Range("D15").Value = "'66.51"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "'0.000007965"
$ws.Range("E16").Value = "  -4.86%  "
$ws.Range("D17").Value = "26.151.52"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "'1.003"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'4.713"
$ws.Range("E19").Value = "  -3.82%  "
$ws.Range("D20").Value = "'186.74"
$ws.Range("E20").Value = "  -3.48%  "
$ws.Range("D21").Value = "'10.26"
$ws.Range("E21").Value = "  -5.42%  "
$ws.Range("D22").Value = "'6.192"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "'149.56"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "'0.1246"
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("D26").Value = "'7.493"
$ws.Range("E26").Value = "  -4.56%  "
$ws.Range("D27").Value = "'15.87"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "'0.06371"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").Value = "'1.351"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("D31").Value = "'3.488"
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("D32").Value = "'3.416"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("D33").Value = "'1.638"
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("D34").Value = "'1.004"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("D35").Value = "'2.409"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "'0.6006"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("D37").Value = "'2.745"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "1.110.21"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'6.111"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").Value = "'0.01614"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").Value = "'0.8656"
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("D42").Value = "'1.004"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").Value = "'99.95"
$ws.Range("D44").Value = "1.817.89"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = "'0.00000000109"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("D46").Value = "'55.26"
$ws.Range("E46").Value = "  -4.38%  "
$ws.Range("D47").Value = "'0.9992"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "'8.066"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "'0.05233"
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").Value = "'0.4248"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").Value = "'5.918"
$ws.Range("E51").Value = "  -2.47%  "
